$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(793, 1).Value = "TestHeader"
$ws.Cells.Item(794, 1).Value = "TestAndaman"
$ws.Cells.Item(794, 2).Value = 190
